# Fix Training Data Issue (#48)
# The "Date" column (BF) held a malformed label "6-14-2007-08" on every
# data row. Correct it to the proper ISO date "2008-06-14" for all rows.
#
# NOTE: assigning a date-shaped string straight to Range.Value causes Excel
# to auto-convert it into a serial date number (and apply a date number
# format). To keep the cell a genuine text value with its original
# (default) style, we temporarily force the target range to Text format,
# write the strings, then clear the temporary formatting back off again.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldDate = "6-14-2007-08"
$newDate = "2008-06-14"

$dateRange = $ws.Range("BF2:BF31")

# Force text interpretation so the new value isn't reinterpreted as a date.
$dateRange.NumberFormat = "@"

for ($row = 2; $row -le 31; $row++) {
    $cell = $ws.Range("BF$row")
    if ($cell.Value() -eq $oldDate) {
        $cell.Value = $newDate
    }
}

# Remove the temporary text formatting we applied so the cells end up with
# their original (default) style, unchanged apart from the text itself.
$dateRange.ClearFormats()
